$wb = $excel.ActiveWorkbook

# --- Sheet 1: By_Odds_Bin ---
$ws1 = $wb.Worksheets.Item("By_Odds_Bin")

$ws1.Range("B2").Value = 15
$ws1.Range("C2").Value = 0
$ws1.Range("D2").Value = 12
$ws1.Range("E2").Value = -12
$ws1.Range("F2").Value = 20

$ws1.Range("B3").Value = 71
$ws1.Range("C3").Value = 0
$ws1.Range("D3").Value = 60
$ws1.Range("E3").Value = -60
$ws1.Range("F3").Value = 15.5

# --- Sheet 2: By_Field_Size ---
$ws2 = $wb.Worksheets.Item("By_Field_Size")

$ws2.Range("B2").Value = 15
$ws2.Range("C2").Value = -8.5
$ws2.Range("D2").Value = 5.5
$ws2.Range("E2").Value = -14
$ws2.Range("F2").Value = 6.7

$ws2.Range("B3").Value = 20
$ws2.Range("C3").Value = 9.5
$ws2.Range("D3").Value = 24.5
$ws2.Range("E3").Value = -15
$ws2.Range("F3").Value = 25

$ws2.Range("B4").Value = 25
$ws2.Range("C4").Value = 8
$ws2.Range("D4").Value = 28
$ws2.Range("E4").Value = -20
$ws2.Range("F4").Value = 20

$ws2.Range("B5").Value = 26
$ws2.Range("C5").Value = -9
$ws2.Range("D5").Value = 14
$ws2.Range("E5").Value = -23
$ws2.Range("F5").Value = 11.5

$ws2.Range("B6").Value = 0
$ws2.Range("C6").Value = 0
$ws2.Range("D6").Value = 0
$ws2.Range("E6").Value = 0
$ws2.Range("F6").Value = ""

# --- Sheet 3: By_Track ---
$ws3 = $wb.Worksheets.Item("By_Track")

$ws3.Range("A2").Value = "NEWTON ABBOT"
$ws3.Range("B2").Value = 86
$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 72
$ws3.Range("E2").Value = -72
$ws3.Range("F2").Value = 16.3
